# Apply marksheet correction: update Right-marks and Total/Max values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right column (B11) 3 -> 5
$ws.Cells.Item(11, 2).Value = 5

# Row 12 "Total": Right column (B12) 42 -> 70
$ws.Cells.Item(12, 2).Value = 70

# Row 12 "Total": Max column (E12) "40/84" -> "70/140"
$ws.Cells.Item(12, 5).Value = "70/140"
